# bad-headers.xlsx fixture update
#
# 1. Header cell A1 on "Main root" changes from "Id" to "Identifier"
#    (a new shared string "Identifier" is introduced ahead of "Id" in the
#    canonical shared-string table, which is what the cell picks up).
# 2. The "Main root" sheet view's bottom-right pane selection moves from
#    A6 back to A2.
# 3. Each of the four worksheets gets an extra (duplicate) filter-database
#    defined name "_xlnm._FilterDatabase_0_0", mirroring the existing
#    "_xlnm._FilterDatabase" / "_xlnm._FilterDatabase_0" pair already
#    present for every sheet.

$wb = $excel.ActiveWorkbook

$wsMainRoot = $wb.Worksheets.Item("Main root")
$wsNodes = $wb.Worksheets.Item("Nodes")
$wsLeaves = $wb.Worksheets.Item("Leaves")
$wsOneToMany = $wb.Worksheets.Item("One to many rows")

# 1. Rename the "Id" header on the "Main root" sheet to "Identifier".
$wsMainRoot.Range("A1").Value = "Identifier"

# 2. Restore the bottom-right pane's selection to A2.
[void]$wsMainRoot.Activate()
[void]$wsMainRoot.Range("A2").Select()

# 3. Re-create a filter database definition for every sheet, which Excel
#    names "_xlnm._FilterDatabase_0_0" since "_xlnm._FilterDatabase" and
#    "_xlnm._FilterDatabase_0" already exist for each of them.
[void]$wsMainRoot.Names.Add("_xlnm._FilterDatabase_0_0", "='Main root'!`$A`$1:`$B`$2")
[void]$wsNodes.Names.Add("_xlnm._FilterDatabase_0_0", "=Nodes!`$A`$1:`$D`$4")
[void]$wsLeaves.Names.Add("_xlnm._FilterDatabase_0_0", "=Leaves!`$A`$1:`$F`$7")
[void]$wsOneToMany.Names.Add("_xlnm._FilterDatabase_0_0", "='One to many rows'!`$A`$1:`$A`$13")
